$wb = $excel.ActiveWorkbook

function Set-CellValue($ws, [string]$cellRef, $newValue) {
    $ws.Range($cellRef).Value = $newValue
}

# --- Worksheet: Citywide Totals ---
$ws = $wb.Worksheets.Item('Citywide Totals')
Set-CellValue $ws "J2" 3260   # was 3236
Set-CellValue $ws "J3" 3397   # was 3377
Set-CellValue $ws "E4" 1991   # was 1990
Set-CellValue $ws "J4" 752   # was 749
Set-CellValue $ws "J6" 4009   # was 3984
Set-CellValue $ws "E7" 25995   # was 25994
Set-CellValue $ws "J7" 11681   # was 11609

# --- Worksheet: Bridgeport ---
$ws = $wb.Worksheets.Item('Bridgeport')
Set-CellValue $ws "J2" 20   # was 19
Set-CellValue $ws "J7" 49   # was 48

# --- Worksheet: Grand Crossing ---
$ws = $wb.Worksheets.Item('Grand Crossing')
Set-CellValue $ws "J2" 111   # was 110
Set-CellValue $ws "J3" 131   # was 130
Set-CellValue $ws "J6" 105   # was 106
Set-CellValue $ws "J7" 371   # was 370

# --- Worksheet: Woodlawn ---
$ws = $wb.Worksheets.Item('Woodlawn')
Set-CellValue $ws "J6" 47   # was 46
Set-CellValue $ws "J7" 168   # was 167

# --- Worksheet: North Lawndale ---
$ws = $wb.Worksheets.Item('North Lawndale')
Set-CellValue $ws "J2" 88   # was 87
Set-CellValue $ws "J3" 171   # was 170
Set-CellValue $ws "J6" 119   # was 117
Set-CellValue $ws "J7" 417   # was 413

# --- Worksheet: New City ---
$ws = $wb.Worksheets.Item('New City')
Set-CellValue $ws "J4" 13   # was 14
Set-CellValue $ws "J7" 306   # was 307

# --- Worksheet: By Neighborhood ---
$ws = $wb.Worksheets.Item('By Neighborhood')
Set-CellValue $ws "J7" 356   # was 353
Set-CellValue $ws "J8" 753   # was 746
Set-CellValue $ws "J11" 172   # was 169
Set-CellValue $ws "J14" 49   # was 48
Set-CellValue $ws "J15" 135   # was 134
Set-CellValue $ws "J19" 363   # was 361
Set-CellValue $ws "J20" 247   # was 245
Set-CellValue $ws "J27" 70   # was 69
Set-CellValue $ws "J29" 676   # was 672
Set-CellValue $ws "J33" 512   # was 505
Set-CellValue $ws "J37" 371   # was 370
Set-CellValue $ws "J41" 77   # was 76
Set-CellValue $ws "J42" 465   # was 462
Set-CellValue $ws "J48" 115   # was 114
Set-CellValue $ws "J50" 65   # was 64
Set-CellValue $ws "J52" 322   # was 319
Set-CellValue $ws "J54" 223   # was 221
Set-CellValue $ws "E63" 335   # was 334
Set-CellValue $ws "J63" 56   # was 51
Set-CellValue $ws "J65" 306   # was 307
Set-CellValue $ws "J67" 417   # was 413
Set-CellValue $ws "J71" 43   # was 42
Set-CellValue $ws "J77" 98   # was 97
Set-CellValue $ws "J79" 348   # was 347
Set-CellValue $ws "J83" 271   # was 268
Set-CellValue $ws "J85" 533   # was 529
Set-CellValue $ws "J90" 137   # was 135
Set-CellValue $ws "J91" 133   # was 131
Set-CellValue $ws "J94" 102   # was 103
Set-CellValue $ws "J95" 182   # was 181
Set-CellValue $ws "J97" 72   # was 68
Set-CellValue $ws "J98" 70   # was 68
Set-CellValue $ws "J99" 168   # was 167
Set-CellValue $ws "E101" 25995   # was 25994
Set-CellValue $ws "J101" 11681   # was 11609

# --- Worksheet: South Chicago ---
$ws = $wb.Worksheets.Item('South Chicago')
Set-CellValue $ws "J3" 102   # was 101
Set-CellValue $ws "J6" 75   # was 73
Set-CellValue $ws "J7" 271   # was 268

# --- Worksheet: West Pullman ---
$ws = $wb.Worksheets.Item('West Pullman')
Set-CellValue $ws "J3" 55   # was 54
Set-CellValue $ws "J7" 182   # was 181

# --- Worksheet: Garfield Park ---
$ws = $wb.Worksheets.Item('Garfield Park')
Set-CellValue $ws "J2" 138   # was 135
Set-CellValue $ws "J3" 165   # was 162
Set-CellValue $ws "J6" 166   # was 165
Set-CellValue $ws "J7" 512   # was 505

# --- Worksheet: Loop ---
$ws = $wb.Worksheets.Item('Loop')
Set-CellValue $ws "J2" 58   # was 56
Set-CellValue $ws "J7" 223   # was 221

# --- Worksheet: Englewood ---
$ws = $wb.Worksheets.Item('Englewood')
Set-CellValue $ws "J2" 207   # was 205
Set-CellValue $ws "J3" 232   # was 231
Set-CellValue $ws "J6" 170   # was 169
Set-CellValue $ws "J7" 676   # was 672

# --- Worksheet: Chatham ---
$ws = $wb.Worksheets.Item('Chatham')
Set-CellValue $ws "J2" 85   # was 84
Set-CellValue $ws "J6" 140   # was 139
Set-CellValue $ws "J7" 363   # was 361

# --- Worksheet: Lake View ---
$ws = $wb.Worksheets.Item('Lake View')
Set-CellValue $ws "J2" 22   # was 21
Set-CellValue $ws "J7" 115   # was 114

# --- Worksheet: South Shore ---
$ws = $wb.Worksheets.Item('South Shore')
Set-CellValue $ws "J2" 131   # was 130
Set-CellValue $ws "J6" 151   # was 148
Set-CellValue $ws "J7" 533   # was 529

# --- Worksheet: Hermosa ---
$ws = $wb.Worksheets.Item('Hermosa')
Set-CellValue $ws "J3" 13   # was 12
Set-CellValue $ws "J7" 77   # was 76

# --- Worksheet: Humboldt Park ---
$ws = $wb.Worksheets.Item('Humboldt Park')
Set-CellValue $ws "J3" 104   # was 103
Set-CellValue $ws "J4" 23   # was 22
Set-CellValue $ws "J6" 230   # was 229
Set-CellValue $ws "J7" 465   # was 462

# --- Worksheet: Washington Park ---
$ws = $wb.Worksheets.Item('Washington Park')
Set-CellValue $ws "J3" 62   # was 61
Set-CellValue $ws "J6" 21   # was 20
Set-CellValue $ws "J7" 133   # was 131

# --- Worksheet: Roseland ---
$ws = $wb.Worksheets.Item('Roseland')
Set-CellValue $ws "J3" 128   # was 127
Set-CellValue $ws "J7" 348   # was 347

# --- Worksheet: Chicago Lawn ---
$ws = $wb.Worksheets.Item('Chicago Lawn')
Set-CellValue $ws "J2" 80   # was 79
Set-CellValue $ws "J3" 77   # was 76
Set-CellValue $ws "J7" 247   # was 245

# --- Worksheet: Little Village ---
$ws = $wb.Worksheets.Item('Little Village')
Set-CellValue $ws "J3" 93   # was 90
Set-CellValue $ws "J7" 322   # was 319

# --- Worksheet: West Loop ---
$ws = $wb.Worksheets.Item('West Loop')
Set-CellValue $ws "J6" 53   # was 54
Set-CellValue $ws "J7" 102   # was 103

# --- Worksheet: Brighton Park ---
$ws = $wb.Worksheets.Item('Brighton Park')
Set-CellValue $ws "J3" 37   # was 36
Set-CellValue $ws "J7" 135   # was 134

# --- Worksheet: Wicker Park ---
$ws = $wb.Worksheets.Item('Wicker Park')
Set-CellValue $ws "J3" 10   # was 9
Set-CellValue $ws "J6" 37   # was 36
Set-CellValue $ws "J7" 70   # was 68

# --- Worksheet: Lincoln Square ---
$ws = $wb.Worksheets.Item('Lincoln Square')
Set-CellValue $ws "J2" 18   # was 17
Set-CellValue $ws "J7" 65   # was 64

# --- Worksheet: Belmont Cragin ---
$ws = $wb.Worksheets.Item('Belmont Cragin')
Set-CellValue $ws "J2" 65   # was 62
Set-CellValue $ws "J7" 172   # was 169

# --- Worksheet: West Town ---
$ws = $wb.Worksheets.Item('West Town')
Set-CellValue $ws "J6" 42   # was 38
Set-CellValue $ws "J7" 72   # was 68

# --- Worksheet: Austin ---
$ws = $wb.Worksheets.Item('Austin')
Set-CellValue $ws "J2" 224   # was 223
Set-CellValue $ws "J3" 236   # was 235
Set-CellValue $ws "J6" 230   # was 225
Set-CellValue $ws "J7" 753   # was 746

# --- Worksheet: Edgewater ---
$ws = $wb.Worksheets.Item('Edgewater')
Set-CellValue $ws "J6" 26   # was 25
Set-CellValue $ws "J7" 70   # was 69

# --- Worksheet: Washington Heights ---
$ws = $wb.Worksheets.Item('Washington Heights')
Set-CellValue $ws "J2" 44   # was 43
Set-CellValue $ws "J6" 42   # was 41
Set-CellValue $ws "J7" 137   # was 135

# --- Worksheet: Oakland ---
$ws = $wb.Worksheets.Item('Oakland')
Set-CellValue $ws "J4" 2   # was 1
Set-CellValue $ws "J7" 43   # was 42

# --- Worksheet: Riverdale ---
$ws = $wb.Worksheets.Item('Riverdale')
Set-CellValue $ws "J3" 33   # was 32
Set-CellValue $ws "J7" 98   # was 97

# --- Worksheet: Auburn Gresham ---
$ws = $wb.Worksheets.Item('Auburn Gresham')
Set-CellValue $ws "J2" 122   # was 120
Set-CellValue $ws "J6" 115   # was 114
Set-CellValue $ws "J7" 356   # was 353
